$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, same style as the other header cells (bold, centered, bordered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save flag values for each data row
$saveValues = @(0, 1, 0, 1, 1, 1, 1, 0, 1, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
